# Update column G ("K" = strikeouts) with newly regenerated values,
# replacing the previous "Strike#" derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 3
    4  = 4
    5  = 1
    6  = 2
    7  = 2
    8  = 2
    9  = 1
    10 = 0
    12 = 3
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
